$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (D, E, H) changes. Only rows whose numeric columns change are listed;
# all rows from 3 to 21 also get their date separator changed from "/" to "-".
$rowChanges = @{
    3  = @{ D = 1; G = 1 }
    4  = @{ D = 1; E = 1; H = 0 }
    6  = @{ D = 1; E = 1; H = 0 }
    7  = @{ D = 1; E = 1; H = 0 }
    10 = @{ D = 1; E = 1; H = 0 }
    11 = @{ D = 1; E = 1; H = 0 }
    13 = @{ D = 1; E = 1; H = 0 }
}

for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Range("A$r")
    $oldDate = $cell.Value2
    $newDate = $oldDate -replace "/", "-"
    # Force text so Excel doesn't reinterpret the dd-mm-yyyy string as a date serial.
    $cell.NumberFormat = "@"
    $cell.Value = $newDate

    if ($rowChanges.ContainsKey($r)) {
        $changes = $rowChanges[$r]
        foreach ($col in $changes.Keys) {
            $ws.Range("$col$r").Value = $changes[$col]
        }
    }
}
